$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Probe1")
$ws2 = $wb.Worksheets.Item("Probe2")

# Zugriffsbeschränkung values: the placeholder test strings are replaced
# with the real allowed values for this field.
$ws1.Range("B34").Value = "restriktiv"
$ws2.Range("B34").Value = "offen"

# Leave the selection on B34 (the just-edited cell) on both sheets, and
# make "Probe2" the active tab when the workbook is reopened.
$ws1.Activate()
$ws1.Range("B34").Select()

$ws2.Activate()
$ws2.Range("B34").Select()
